$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.404.89"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "2.306.36"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "311.24"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "103.21"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +8.48%  "
$ws.Range("D10").Value = "36.55"
$ws.Range("E10").Value = "  +4.48%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("D12").Value = "51.90"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "7.05"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "2.664.01"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "15.12"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("D17").Value = "2.311.82"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "0.811"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("D19").Value = "43.293.72"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "12.25"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("E22").Value = "  +4.27%  "
$ws.Range("D23").Value = "68.14"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "243.24"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("E25").Value = "  +3.51%  "
$ws.Range("D26").Value = "2.61"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "24.78"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("E29").Value = "  +8.48%  "
$ws.Range("D30").Value = "36.98"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").Value = "168.47"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "18.10"
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("D36").Value = "2.53"
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").Value = "0.0743"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "3.05"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  +4.02%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "4.48"
$ws.Range("E40").Value = "  +9.33%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.106"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  +13.89%  "
$ws.Range("D44").Value = "0.0294"
$ws.Range("E44").Value = "  +5.54%  "
$ws.Range("D45").Value = "1.985.36"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "19.00"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").Value = "3.01"
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").Value = "10.02"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "56.02"
$ws.Range("E49").Value = "  +5.51%  "
$ws.Range("D50").Value = "2.94"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  +9.43%  "
